$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("A_SERIES_ERSTAUSB", "K_QUALI", "Erstausbildung", "XXX")
  ,@("A_SERIES_FACHKR", "K_SERIES", "Fachkräftepotenzial (insgesamt)", "XXX")
  ,@("A_SERIES_FINENTRANS", "K_SERIES", "Endenergieverbrauch durch Güterbeförderung", "Final energy consumption for the transport of goods")
  ,@("A_SERIES_FOREIGN", "K_SERIES", "Ausländer", "Foreigns")
  ,@("A_SERIES_FTTB", "K_SERIES", "Reine Glasfasernetze (FTTB/H)", "Fibre-to-the-Building/Home (FTTB/H)")
  ,@("A_SERIES_GDP", "K_SERIES", "Bruttoinlandsprodukt (preisbereinigt)", "Gross domestic product (price-adjusted)")
  ,@("A_SERIES_GERMAN", "K_SERIES", "Deutsche", "Germans")
  ,@("A_SERIES_GOODSTRANS", "K_SERIES", "Güterbeförderungsleistung", "Goods transport performance")
  ,@("A_SERIES_HOCHSCHUL", "K_QUALI", "Hochschulabschlüsse", "XXX")
  ,@("A_SERIES_HOEHERQUAL", "K_QUALI", "Höherqualifizierende Berufsausbildung", "XXX")
  ,@("A_SERIES_IMP", "K_SERIES", "Einfuhren", "Imports")
  ,@("A_SERIES_LOCATIONS", "K_SERIES", "EMAS-Organisationsstandorte", "Locations of organisation registered with EMAS")
  ,@("A_SERIES_MARKETINC", "K_SERIES", "Gini-Koeffizient des Markteinkommens", "Gini coefficient of market income")
  ,@("A_SERIES_MATINPTCON", "K_SERIES", "Rohstoffeinsatz für Konsum, Investitionen und Exporte", "Raw material input for consumption, investment and exports")
  ,@("A_SERIES_MATINPTPROD", "K_SERIES", "Gesamtrohstoffproduktivität", "Raw material input productivity")
  ,@("A_SERIES_MSYEXAMINED", "K_SERIES", "MSY-untersuchte Fischbestände an allen bewirtschafteten Beständen", "Share of MSY examined stocks in all managed fish stocks")
  ,@("A_SERIES_OBESE", "K_SERIES", "Adipositas", "Obese")
  ,@("A_SERIES_ORGANISAT", "K_SERIES", "EMAS-Organisationen", "EMAS organisations")
  ,@("A_SERIES_OVERWEIGHT", "K_SERIES", "Übergewicht ohne Adipositas", "Overweight without obese")
  ,@("A_SERIES_PANDEMIC", "K_SERIES", "Pandemieprävention und -reaktion", "Pandemic prevention and response")
  ,@("A_SERIES_PASSTRANS", "K_SERIES", "Personenbeförderungsleistung", "Passenger transport performance")
  ,@("A_SERIES_POSTSEC", "K_SERIES", "Postsekundarer nicht-tertiärer Abschluss", "With post-secondary non-tertiary qualification")
  ,@("A_SERIES_PREV", "K_SERIES", "Durch Prävention vermeidbaren Sterblichkeit", "Mortality avoidable through prevention")
  ,@("A_SERIES_PREVUNDBEH", "K_SERIES", "Durch Prävention und Behandlung vermeidbaren Sterblichkeit", "Mortality avoidable through prevention and treatment")
  ,@("A_SERIES_PROCIMP", "K_SERIES", "Einfuhren von weiterverarbeiteten Produkten", "Imports of processed products")
  ,@("A_SERIES_QUALIFZIERT", "K_SERIES", "Akademisch Qualifizierte und beruflich Höherqualifizierte", "XXX")
  ,@("A_SERIES_RESLDC", "K_SERIES", "Forschende aus LDCs", "Researchers from LDCs")
  ,@("A_SERIES_RESOTHER", "K_SERIES", "Forschende aus den übrigen Entwicklungs- und Schwellenländern", "Researchers from other developing and emerging countries")
  ,@("A_SERIES_STUDLDC", "K_SERIES", "Studierende aus LDCs", "Students from LDCs")
  ,@("A_SERIES_STUDOTHER", "K_SERIES", "Studierende aus den übrigen Entwicklungs- und Schwellenländern", "Students from other developing and emerging countries")
  ,@("A_SERIES_TERTIARY", "K_SERIES", "Tertiärer Abschluss", "With tertiary qualification")
  ,@("A_SERIES_TOTAL", "K_SERIES", "Insgesamt", "Total")
  ,@("A_SERIES_TOTALOBJ", "K_SERIES", "Objekte insgesamt", "All objects")
  ,@("A_SERIES_VALUECONSUM", "K_SERIES", "Wert von Konsum, Investitionen und Exporten (preisbereinigt)", "Value of consumption, investment and exports (price-adjusted)")
  ,@("A_SERIES_WEALTH", "K_SERIES", "Gini-Koeffizient des Vermögens", "Gini coefficient of wealth")
  ,@("A_SERIES_WITHOUTDIGI", "K_SERIES", "Objekte ohne Digitalisat", "Objects without digitised media")
  ,@("A_SERIES_YEAR2012", "K_SERIES", "Bezogen auf die Mittel- und Oberzentren des Jahres 2012", "For the medium-sized and major cities of 2012")
  ,@("A_SERIES_YEARINQUEST", "K_SERIES", "Bezogen auf die Mittel- und Oberzentren des jeweiligen Jahres", "For the medium-sized and major cities of the year in question")
  ,@("A_SERIES_ZWEITAUSB", "K_QUALI", "Zweitausbildung, Ausbildung nach Abitur", "XXX")
  ,@("A_SEX_D", "K_SEX", "Divers", "Divers")
  ,@("A_SEX_F", "K_SEX", "Weiblich", "Female")
  ,@("A_SEX_M", "K_SEX", "Männlich", "Male")
  ,@("A_SEX_U", "K_SEX", "Unbekannt", "Unknown")
  ,@("A_SUBINDEX_COAST", "K_SUBINDEX", "Teilindex Küsten/Meere", "Coasts and seas sub-index")
  ,@("A_SUBINDEX_FARMLAND", "K_SUBINDEX", "Teilindex Agrarland", "Farmland sub-index")
  ,@("A_SUBINDEX_FOREST", "K_SUBINDEX", "Teilindex Wälder", "Forest sub-index")
  ,@("A_SUBINDEX_SETTLEMNT", "K_SUBINDEX", "Teilindex Siedlungen", "Settlements sub-index")
  ,@("A_SUBINDEX_WATERS", "K_SUBINDEX", "Teilindex Binnengewässer", "Inland waters sub-index")
  ,@("A_TYPEAREA_RECREATIO", "K_TYPEAREA", "Erholungsfläche, Friedhof", "Recreation area, cemetery")
  ,@("A_TYPEAREA_SETTLEMEN", "K_TYPEAREA", "Siedlungsfläche (bebaut)", "Settlement area (built on)")
  ,@("A_TYPEAREA_TRANSPORT", "K_TYPEAREA", "Verkehrsfläche", "Transport area")
  ,@("A_URBAN_NONRURAL", "K_URBAN", "Nicht-ländliche Gebiete", "Non-rural areas")
  ,@("A_URBAN_RURAL", "K_URBAN", "Ländliche Gebiete", "Rural areas")
)

$startRow = 483
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $rec = $data[$i]
  $ws.Cells.Item($row, 1).Value = $rec[0]
  $ws.Cells.Item($row, 2).Value = $rec[1]
  $ws.Cells.Item($row, 3).Value = $rec[2]
  $ws.Cells.Item($row, 4).Value = $rec[3]
}
